$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two rows (18 and 19) so the roster shrinks from 18 to 16 players
$ws.Range("A18:C19").Delete("xlShiftUp")

# Update the remaining 16 data rows (2-17) with the new roster
$ws.Range("A2").Value2 = "De'Aaron Fox"
$ws.Range("B2").Value2 = "PG"
$ws.Range("C2").Value2 = "Sacramento Kings"

$ws.Range("A3").Value2 = "Shaedon Sharpe"
$ws.Range("B3").Value2 = "SG,SF"
$ws.Range("C3").Value2 = "Portland Trail Blazers"

$ws.Range("A4").Value2 = "Josh Giddey"
$ws.Range("B4").Value2 = "PG,SG,SF"
$ws.Range("C4").Value2 = "Chicago Bulls"

$ws.Range("A5").Value2 = "Tyler Herro"
$ws.Range("B5").Value2 = "PG,SG"
$ws.Range("C5").Value2 = "Miami Heat"

$ws.Range("A6").Value2 = "DeMar DeRozan"
$ws.Range("B6").Value2 = "SF,PF"
$ws.Range("C6").Value2 = "Sacramento Kings"

$ws.Range("A7").Value2 = "Kyle Kuzma"
$ws.Range("B7").Value2 = "PF"
$ws.Range("C7").Value2 = "Washington Wizards"

$ws.Range("A8").Value2 = "Luguentz Dort"
$ws.Range("B8").Value2 = "SG,SF"
$ws.Range("C8").Value2 = "Oklahoma City Thunder"

$ws.Range("A9").Value2 = "Nikola Vucevic"
$ws.Range("B9").Value2 = "PF,C"
$ws.Range("C9").Value2 = "Chicago Bulls"

$ws.Range("A10").Value2 = "Brook Lopez"
$ws.Range("B10").Value2 = "C"
$ws.Range("C10").Value2 = "Milwaukee Bucks"

$ws.Range("A11").Value2 = "Nick Richards"
$ws.Range("B11").Value2 = "C"
$ws.Range("C11").Value2 = "Phoenix Suns"

$ws.Range("A12").Value2 = "Miles Bridges"
$ws.Range("B12").Value2 = "SF,PF"
$ws.Range("C12").Value2 = "Charlotte Hornets"

$ws.Range("A13").Value2 = "Evan Mobley"
$ws.Range("B13").Value2 = "PF,C"
$ws.Range("C13").Value2 = "Cleveland Cavaliers"

$ws.Range("A14").Value2 = "Scottie Barnes"
$ws.Range("B14").Value2 = "PG,SG,SF,PF"
$ws.Range("C14").Value2 = "Toronto Raptors"

$ws.Range("A15").Value2 = "Mikal Bridges"
$ws.Range("B15").Value2 = "SG,SF,PF"
$ws.Range("C15").Value2 = "New York Knicks"

$ws.Range("A16").Value2 = "Luka Doncic"
$ws.Range("B16").Value2 = "PG,SG"
$ws.Range("C16").Value2 = "Dallas Mavericks"

$ws.Range("A17").Value2 = "Ja Morant"
$ws.Range("B17").Value2 = "PG"
$ws.Range("C17").Value2 = "Memphis Grizzlies"
